$wb = $excel.ActiveWorkbook

# Sheet "new_users": Interest (D2) changes from ManualTesting to Teamcity,
# and the active selection moves to D2.
$wsUsers = $wb.Worksheets.Item("new_users")
$wsUsers.Range("D2").Value = "Teamcity"
$wsUsers.Range("D2").Select()

# Sheet "addCourse_details": CategoryName (H2) changes from Informatica to
# Teamcity, and the active selection moves to H9.
$wsCourse = $wb.Worksheets.Item("addCourse_details")
$wsCourse.Range("H2").Value = "Teamcity"
$wsCourse.Range("H9").Select()
